$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Existing": mark FEINSSN present (F17) for the "Sales and Use Tax"
# row that was just added, and leave the cursor on D23 (next empty row).
# ---------------------------------------------------------------------------
$wsExisting = $wb.Worksheets.Item("Existing")
$wsExisting.Range("F17").Value = "Y"

# ---------------------------------------------------------------------------
# Sheet "NewTaxReturn": record the results of three additional RAD test runs
# (rows 27-29) and flag the FEINSSN column for a few rows that already run.
# ---------------------------------------------------------------------------
$wsNewTaxReturn = $wb.Worksheets.Item("NewTaxReturn")

$wsNewTaxReturn.Range("G26").Value = "Y"

$wsNewTaxReturn.Range("A27").Value = "Pass"
$wsNewTaxReturn.Range("B27").Value = "Mon Jan 15 21:03:38 EST 2024"
$wsNewTaxReturn.Range("A27").Style = "Normal"
$wsNewTaxReturn.Range("B27").Style = "Normal"

$wsNewTaxReturn.Range("A28").Value = "Pass"
$wsNewTaxReturn.Range("B28").Value = "Mon Jan 15 21:05:26 EST 2024"
$wsNewTaxReturn.Range("A28").Style = "Normal"
$wsNewTaxReturn.Range("B28").Style = "Normal"

$wsNewTaxReturn.Range("A29").Value = "Pass"
$wsNewTaxReturn.Range("B29").Value = "Mon Jan 15 21:07:19 EST 2024"
$wsNewTaxReturn.Range("A29").Style = "Normal"
$wsNewTaxReturn.Range("B29").Style = "Normal"

$wsNewTaxReturn.Range("G39").Value = "Y"
$wsNewTaxReturn.Range("G52").Value = "Y"

# Column B ("Date") now holds longer timestamp text - best-fit the width.
$wsNewTaxReturn.Columns("B").AutoFit()

# Leave the selection on the "Execute" column for the full data range.
$wsNewTaxReturn.Range("C2:C55").Select()

# ---------------------------------------------------------------------------
# Make "Personal_EL" the active / visible tab, as it was the sheet the author
# ended up looking at after running the RAD tests.
# ---------------------------------------------------------------------------
$wsPersonalEL = $wb.Worksheets.Item("Personal_EL")
$wsPersonalEL.Activate()

# Restore the cursor on "Existing" to D23 (recorded in its own sheet view).
$wsExisting.Range("D23").Select()
$wsPersonalEL.Activate()
